# The sheet's weekly "Papa" (potato) price table at
# Macroferia Regional de Talca gained two new daily price entries for the
# most recent date (44706). These are inserted as new rows 466-467, pushing
# every existing row from the old 466 downward by two (466->468, ...,
# 491->493), which is exactly what the diff shows (dimension A1:R491 ->
# A1:R493, with the old row 466.. data now living two rows lower).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the old row 466; Excel shifts every
# row below (466..491) down to (468..493) and the dimension/ref auto-grows.
$ws.Rows("466:467").Insert()

# New row 466: Asterix, "1a (cosecha)"
$ws.Range("A466").Value = 5
$ws.Range("B466").Value = "Macroferia Regional de Talca"
$ws.Range("C466").Value = "Maule"
$ws.Range("D466").Value = 44706
$ws.Range("E466").Value = 7
$ws.Range("F466").Value = 100114001
$ws.Range("G466").Value = "Papa"
$ws.Range("H466").Value = "Asterix"
$ws.Range("I466").Value = "1a (cosecha)"
$ws.Range("J466").Value = 1200
$ws.Range("K466").Value = 7000
$ws.Range("L466").Value = 7000
$ws.Range("M466").Value = 7000
$ws.Range("N466").Value = "$/saco 25 kilos"
$ws.Range("O466").Value = "Región de Los Lagos"
$ws.Range("P466").Value = 280
$ws.Range("Q466").Value = 25
$ws.Range("R466").Value = "Hortaliza"

# New row 467: Rodeo, "1a (cosecha lavada)"
$ws.Range("A467").Value = 5
$ws.Range("B467").Value = "Macroferia Regional de Talca"
$ws.Range("C467").Value = "Maule"
$ws.Range("D467").Value = 44706
$ws.Range("E467").Value = 7
$ws.Range("F467").Value = 100114001
$ws.Range("G467").Value = "Papa"
$ws.Range("H467").Value = "Rodeo"
$ws.Range("I467").Value = "1a (cosecha lavada)"
$ws.Range("J467").Value = 1600
$ws.Range("K467").Value = 8000
$ws.Range("L467").Value = 8000
$ws.Range("M467").Value = 8000
$ws.Range("N467").Value = "$/malla 25 kilos"
$ws.Range("O467").Value = "Región de Los Lagos"
$ws.Range("P467").Value = 320
$ws.Range("Q467").Value = 25
$ws.Range("R467").Value = "Hortaliza"
